# Updates cryptos list values per the Oct 15 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.302.38"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.566.64"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'210.42"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'22.08"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.249"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'0.0871"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "1.789.58"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "1.564.09"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'3.78"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "27.268.02"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "'61.91"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'217.43"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "'9.22"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'152.95"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'15.02"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'0.107"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D34").Value = "1.435.05"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'0.533"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'5.93"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "'2.34"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'64.54"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "1.703.13"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "'0.0955"
$ws.Range("E51").Value = "  -0.47%  "
